# AMOS B07 - Demo Day Preparation: update the table style (theme) applied
# to the two data tables in the deck.
#
# Slide 6  (table "Google Shape;70;p13")  : Table_0  -> built-in style {66B0EB45-EE49-450D-A477-25FF1C32F1AE}
# Slide 14 (table "Google Shape;128;p21") : Table_1  -> built-in style {F574F4D0-4C87-4EFE-9EAA-9E4A408D93CE}

$p = $ppt.ActivePresentation

# --- Slide 6 table -------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$table6Shape = $slide6.Shapes.Item(3)
if ($table6Shape.HasTable) {
    $table6Shape.Table.ApplyStyle("{66B0EB45-EE49-450D-A477-25FF1C32F1AE}")
}

# --- Slide 14 table --------------------------------------------------------
$slide14 = $p.Slides.Item(14)
$table14Shape = $slide14.Shapes.Item(3)
if ($table14Shape.HasTable) {
    $table14Shape.Table.ApplyStyle("{F574F4D0-4C87-4EFE-9EAA-9E4A408D93CE}")
}
